$d = $word.ActiveDocument
$wordNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# 1. Remove the "Meta description" paragraph that currently sits right after
#    the H1 title ("Play Chicken Party Free - Review of Booming Games Slot").
$metaFind = $d.Content.Find
$metaFind.Text = "Meta description"
$metaFind.Forward = $true
$metaFind.Wrap = 1
if ($metaFind.Execute()) {
    $metaRange = $metaFind.Parent
    $metaRange.Expand(4)  # wdParagraph: grow the found range to the whole paragraph (incl. mark)
    $metaRange.Delete()
}

# 2. At the bottom of the document, the paragraph that used to contain the
#    DALLE feature-image prompt now becomes two paragraphs:
#      - a new bold paragraph repeating the page title
#      - the existing (italic) paragraph, but with the meta-description text
#    Replace the whole trailing paragraph (DALLE prompt) in one shot via XML
#    so the run/paragraph structure matches exactly.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$targetRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$newXml = '<w:p xmlns:w="' + $wordNs + '"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chicken Party Free - Review of Booming Games Slot</w:t></w:r></w:p>' + `
          '<w:p xmlns:w="' + $wordNs + '"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Check out our review of Chicken Party by Booming Games and play for free. Featuring chicken-themed design and bonus mode.</w:t></w:r></w:p>'

$targetRange.InsertXML($newXml)
